$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings are preserved as text (matches original inlineStr cells)
$ws.Range("D2:D50").NumberFormat = "@"
$ws.Range("E2:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Column G ("Hora") moves uniformly from 7 to 8 for every data row
$ws.Range("G2:G51").Value = "8"

# Apply new per-row Price (D) / Volume(1h) (E) values
$ws.Range("D2").Value = "300.06"
$ws.Range("D3").Value = "31.66"
$ws.Range("D4").Value = "5.129"
$ws.Range("D5").Value = "0.07876"
$ws.Range("D6").Value = "2.270"
$ws.Range("D7").Value = "7.919"
$ws.Range("D8").Value = "3.851"
$ws.Range("D9").Value = "0.9108"
$ws.Range("D10").Value = "0.1740"
$ws.Range("D11").Value = "0.07311"
$ws.Range("D12").Value = "0.08067"
$ws.Range("D13").Value = "0.03040"
$ws.Range("D14").Value = "0.09856"
$ws.Range("D15").Value = "0.001509"
$ws.Range("D16").Value = "0.006077"
$ws.Range("D17").Value = "3.478"
$ws.Range("D18").Value = "2.238"
$ws.Range("D19").Value = "0.3252"
$ws.Range("D20").Value = "0.1350"
$ws.Range("D21").Value = "4.697"
$ws.Range("D22").Value = "0.1657"
$ws.Range("D23").Value = "0.04622"
$ws.Range("D24").Value = "0.001269"
$ws.Range("D25").Value = "0.004482"
$ws.Range("D26").Value = "0.0001189"
$ws.Range("D27").Value = "0.0003440"
$ws.Range("D39").Value = "0.01867"
$ws.Range("D40").Value = "0.04560"
$ws.Range("D41").Value = "0.007354"
$ws.Range("D42").Value = "0.1341"
$ws.Range("D43").Value = "0.002157"
$ws.Range("D44").Value = "0.01070"
$ws.Range("D45").Value = "0.00006409"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D47").Value = "0.006663"
$ws.Range("D49").Value = "0.00002107"
$ws.Range("D50").Value = "0.0002006"
$ws.Range("E2").Value = "1.05%"
$ws.Range("E3").Value = "0.57%"
$ws.Range("E4").Value = "0.17%"
$ws.Range("E5").Value = "5.30%"
$ws.Range("E6").Value = "34.34%"
$ws.Range("E7").Value = "2.38%"
$ws.Range("E8").Value = "1.53%"
$ws.Range("E9").Value = "-2.18%"
$ws.Range("E10").Value = "2.77%"
$ws.Range("E11").Value = "2.70%"
$ws.Range("E12").Value = "1.16%"
$ws.Range("E13").Value = "0.31%"
$ws.Range("E14").Value = "-0.49%"
$ws.Range("E15").Value = "1.09%"
$ws.Range("E16").Value = "-2.78%"
$ws.Range("E17").Value = "0.63%"
$ws.Range("E18").Value = "0.44%"
$ws.Range("E19").Value = "-0.90%"
$ws.Range("E20").Value = "1.12%"
$ws.Range("E21").Value = "3.00%"
$ws.Range("E22").Value = "6.81%"
$ws.Range("E23").Value = "-0.49%"
$ws.Range("E24").Value = "3.97%"
$ws.Range("E25").Value = "1.22%"
$ws.Range("E26").Value = "-8.60%"
$ws.Range("E27").Value = "83.29%"
$ws.Range("E39").Value = "11.45%"
$ws.Range("E40").Value = "2.48%"
$ws.Range("E41").Value = "3.15%"
$ws.Range("E42").Value = "0.96%"
$ws.Range("E43").Value = "4.60%"
$ws.Range("E44").Value = "-13.11%"
$ws.Range("E45").Value = "7.16%"
$ws.Range("E46").Value = "0.31%"
$ws.Range("E47").Value = "-48.61%"
$ws.Range("E48").Value = "-57.47%"
$ws.Range("E49").Value = "0.31%"
$ws.Range("E50").Value = "0.38%"

# Restore default (Normal) style so no stray text-format style index lingers on these cells
$ws.Range("D2:D50").Style = "Normal"
$ws.Range("E2:E50").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
